$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets   = $wb.Worksheets.Item("Assets")

# --- Settings!B3: give it the value "FPPS/DTOE" and reset its formatting to Normal ---
$wsSettings.Range("B3").Style = "Normal"
$wsSettings.Range("B3").Value = "FPPS/DTOE"

# --- Assets!C2:C11: new column of "FPPS/DTOE" values alongside the existing A/B columns ---
$wsAssets.Range("C2:C11").Value = "FPPS/DTOE"

# --- Selection / active-sheet bookkeeping to match the final saved view state ---
[void]$wsAssets.Activate()
[void]$wsAssets.Range("C2").Select()

[void]$wsSettings.Activate()
[void]$wsSettings.Range("B3").Select()
